# Backlog_4.xlsx edit
# Commit: "Semana 04" is dropped from the shared-strings lookup table and the
# "Semana" column (C) on both sheets now stores the plain week number (4)
# instead of the text label. The active sheet / selection state also moves
# from ITI (last row) back to SPN (top of the list).

$wb = $excel.ActiveWorkbook
$wsSPN = $wb.Worksheets.Item("SPN")
$wsITI = $wb.Worksheets.Item("ITI")

# --- Column C ("Semana"): replace the "Semana 04" text with the literal
# number 4 for every data row on both sheets. ---
for ($r = 2; $r -le 30; $r++) {
    $wsSPN.Cells.Item($r, 3).Value = 4
}

for ($r = 2; $r -le 45; $r++) {
    $wsITI.Cells.Item($r, 3).Value = 4
}

# --- View state: SPN becomes the active/selected sheet with C2:C30 selected;
# ITI keeps a C2:C45 selection but is no longer the active tab. ---
$wsITI.Activate()
$wsITI.Range("C2:C45").Select()

$wsSPN.Activate()
$wsSPN.Range("C2:C30").Select()
